$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new value into F4 (adds "sdl" as a new shared string)
$ws.Range("F4").Value = "sdl"

# Move the selection/active cell to F4, matching the author's final cursor position
$ws.Range("F4").Select()
